$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" value (8/31/2019 -> 11/5/2019)
#    on every Date Placeholder: the slide master and all 11 slide layouts.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "8/31/2019") {
                $tr.Text = "11/5/2019"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 6 "Braille Cell Framming  " -> "Braille Cell Framing  "
#    (also collapses the 3 separate runs back into a single run).
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $sh = $slide6.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "Braille Cell Framming  ") {
            $tr = $sh.TextFrame.TextRange
            $tr.Text = "................................................."
            $sh.TextFrame.TextRange.Text = "Braille Cell Framing  "
        }
    }
}
